$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "75.756.74"
$ws.Range("E2").Value = "  +1.18%  "

$ws.Range("D3").Value = "2.919.27"
$ws.Range("E3").Value = "  +4.37%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'200.01"
$ws.Range("E5").Value = "  +6.51%  "

$ws.Range("D6").Value = "'596.40"
$ws.Range("E6").Value = "  +1.06%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.551"
$ws.Range("E8").Value = "  +1.30%  "

$ws.Range("D9").Value = "'0.195"
$ws.Range("E9").Value = "  +3.27%  "

$ws.Range("D10").Value = "2.921.33"
$ws.Range("E10").Value = "  +4.57%  "

$ws.Range("D11").Value = "'0.441"
$ws.Range("E11").Value = "  +17.71%  "

$ws.Range("E12").Value = "  +0.96%  "

$ws.Range("D13").Value = "'4.92"
$ws.Range("E13").Value = "  +1.64%  "

$ws.Range("D14").Value = "3.457.89"
$ws.Range("E14").Value = "  +4.42%  "

$ws.Range("D15").Value = "'28.08"
$ws.Range("E15").Value = "  +4.86%  "

$ws.Range("D16").Value = "75.686.55"
$ws.Range("E16").Value = "  +1.27%  "

$ws.Range("D17").Value = "'0.0000188"
$ws.Range("E17").Value = "  +1.29%  "

$ws.Range("D18").Value = "2.920.18"
$ws.Range("E18").Value = "  +4.70%  "

$ws.Range("D19").Value = "'13.15"
$ws.Range("E19").Value = "  +7.49%  "

$ws.Range("D20").Value = "'8.71"
$ws.Range("E20").Value = "  -2.81%  "

$ws.Range("D21").Value = "'372.55"
$ws.Range("E21").Value = "  -0.83%  "

$ws.Range("D22").Value = "'2.29"
$ws.Range("E22").Value = "  +0.53%  "

$ws.Range("D23").Value = "'4.31"
$ws.Range("E23").Value = "  +5.71%  "

$ws.Range("D24").Value = "'71.83"
$ws.Range("E24").Value = "  +1.61%  "

$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").Value = "3.086.70"
$ws.Range("E26").Value = "  +4.68%  "

$ws.Range("D27").Value = "'4.29"
$ws.Range("E27").Value = "  +3.89%  "

$ws.Range("D28").Value = "'9.65"
$ws.Range("E28").Value = "  +0.62%  "

$ws.Range("D29").Value = "'0.0000107"
$ws.Range("E29").Value = "  +4.82%  "

$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("E31").Value = "  -0.47%  "

$ws.Range("D32").Value = "'7.84"
$ws.Range("E32").Value = "  +3.25%  "

$ws.Range("D33").Value = "'498.92"
$ws.Range("E33").Value = "  -1.88%  "

$ws.Range("D34").Value = "'1.84"
$ws.Range("E34").Value = "  +3.13%  "

$ws.Range("E35").Value = "  +0.14%  "

$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'163.98"
$ws.Range("E36").Value = "  -0.65%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").Value = "'20.17"
$ws.Range("E37").Value = "  +2.18%  "

$ws.Range("D38").Value = "'0.108"
$ws.Range("E38").Value = "  +25.57%  "

$ws.Range("E39").Value = "  +1.35%  "

$ws.Range("D40").Value = "'0.369"
$ws.Range("E40").Value = "  +8.40%  "

$ws.Range("E41").Value = "  -4.29%  "

$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("D43").Value = "'177.93"
$ws.Range("E43").Value = "  -1.17%  "

$ws.Range("D44").Value = "'4.95"
$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("E45").Value = "  -0.30%  "

$ws.Range("D46").Value = "'40.17"
$ws.Range("E46").Value = "  +0.14%  "

$ws.Range("E47").Value = "  -0.90%  "

$ws.Range("D48").Value = "'2.30"
$ws.Range("E48").Value = "  -0.68%  "

$ws.Range("D49").Value = "'0.575"
$ws.Range("E49").Value = "  +1.46%  "

$ws.Range("D50").Value = "'3.83"
$ws.Range("E50").Value = "  +3.32%  "

$ws.Range("D51").Value = "'22.40"
$ws.Range("E51").Value = "  +8.04%  "
